$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated B (count) and C (gain) values per row, per revised feature-importance data for paper graphs
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 9.8255607085309862
$ws.Range("B3").Value = 63
$ws.Range("C3").Value = 406.62451362557482
$ws.Range("B4").Value = 146
$ws.Range("C4").Value = 3888.873594045363
$ws.Range("B5").Value = 31
$ws.Range("C5").Value = 60.380289722136688
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 181.71878996315351
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 8.7714522182941437
$ws.Range("B8").Value = 160
$ws.Range("C8").Value = 830.03624021523228
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 7.5149674545064036
$ws.Range("B10").Value = 25
$ws.Range("C10").Value = 47.648869268596172
$ws.Range("B11").Value = 33
$ws.Range("C11").Value = 101.10505566733789
$ws.Range("B12").Value = 82
$ws.Range("C12").Value = 656.29551499157617
$ws.Range("B13").Value = 50
$ws.Range("C13").Value = 200.25604570625731
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 10.7715220451355
$ws.Range("B15").Value = 23
$ws.Range("C15").Value = 64.533722168882377
$ws.Range("B16").Value = 95
$ws.Range("C16").Value = 1241.52320741117
$ws.Range("B17").Value = 17
$ws.Range("C17").Value = 15.03722193650901
$ws.Range("B18").Value = 33
$ws.Range("C18").Value = 453.91012918018708
$ws.Range("B19").Value = 65
$ws.Range("C19").Value = 759.32909501437098
$ws.Range("B20").Value = 109
$ws.Range("C20").Value = 2432.3685398956509
$ws.Range("B21").Value = 53
$ws.Range("C21").Value = 352.42950563955952
$ws.Range("B22").Value = 10
$ws.Range("C22").Value = 7.5363933009502944
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("B24").Value = 19
$ws.Range("C24").Value = 40.814324602484703
$ws.Range("B25").Value = 15
$ws.Range("C25").Value = 20.579686932265759
$ws.Range("B26").Value = 34
$ws.Range("C26").Value = 70.481048094225116
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("B28").Value = 61
$ws.Range("C28").Value = 189.23404010385269
$ws.Range("B29").Value = 64
$ws.Range("C29").Value = 10283.85591027141
$ws.Range("B30").Value = 24
$ws.Range("C30").Value = 48.266334716542588
$ws.Range("B31").Value = 32
$ws.Range("C31").Value = 140.51849013566971
$ws.Range("B32").Value = 8
$ws.Range("C32").Value = 11.072166204452509
$ws.Range("B33").Value = 5
$ws.Range("C33").Value = 11.53012016415596
$ws.Range("B34").Value = 25
$ws.Range("C34").Value = 43.488101355731487
$ws.Range("B35").Value = 5
$ws.Range("C35").Value = 19.834000110626221
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("B37").Value = 52
$ws.Range("C37").Value = 588.23096197843552
